# PlayerPerformance_5774.xlsx edit:
#  1. Clear the stray empty "B" cells on rows 3-7 of "ODI Batting" so the
#     cells are absent entirely (they were present-but-empty inlineStr cells).
#  2. Add a new worksheet "ODI Batting Extra" at the end of the workbook with
#     header row + data rows (MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#     PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH).

$wb = $excel.ActiveWorkbook

# --- 1. "ODI Batting": drop the empty B3:B7 cells -------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B3:B7").ClearContents()

# --- 2. Add "ODI Batting Extra" sheet at the end ---------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add($null, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Header row (copy formatting/style from an existing header so it matches
# the bold / centered / bordered "header" style already used elsewhere).
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

$batting.Range("A1:F1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)

# Helper: write a text value into a cell, keeping it as TEXT even when the
# literal looks numeric (e.g. "4406"), and without leaving behind a
# quote-prefix / number-format style change on the cell.
function Set-TextCell($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

# Helper: write an explicit empty-string (present but blank) cell.
function Set-EmptyCell($range) {
    Set-TextCell $range ""
}

# Row 2: MATCH_CODE 4406, no batting position, no 4s/6s/pct, not MOM
Set-TextCell $extra.Range("A2") "4406"
Set-EmptyCell $extra.Range("B2")
Set-EmptyCell $extra.Range("C2")
Set-EmptyCell $extra.Range("D2")
Set-EmptyCell $extra.Range("E2")
Set-TextCell $extra.Range("F2") "NO"

# Row 3: MATCH_CODE 4410, batting position 10
Set-TextCell $extra.Range("A3") "4410"
$extra.Range("B3").Value = 10
Set-EmptyCell $extra.Range("C3")
Set-EmptyCell $extra.Range("D3")
Set-EmptyCell $extra.Range("E3")
Set-TextCell $extra.Range("F3") "NO"

# Row 4: MATCH_CODE 4452, batting position 9
Set-TextCell $extra.Range("A4") "4452"
$extra.Range("B4").Value = 9
Set-EmptyCell $extra.Range("C4")
Set-EmptyCell $extra.Range("D4")
Set-EmptyCell $extra.Range("E4")
Set-TextCell $extra.Range("F4") "NO"

# Row 5: MATCH_CODE 4453, batting position 9
Set-TextCell $extra.Range("A5") "4453"
$extra.Range("B5").Value = 9
Set-EmptyCell $extra.Range("C5")
Set-EmptyCell $extra.Range("D5")
Set-EmptyCell $extra.Range("E5")
Set-TextCell $extra.Range("F5") "NO"

# Row 6: MATCH_CODE 4455, batting position 9
Set-TextCell $extra.Range("A6") "4455"
$extra.Range("B6").Value = 9
Set-EmptyCell $extra.Range("C6")
Set-EmptyCell $extra.Range("D6")
Set-EmptyCell $extra.Range("E6")
Set-TextCell $extra.Range("F6") "NO"

# Row 7: MATCH_CODE 4563, batting position 9
Set-TextCell $extra.Range("A7") "4563"
$extra.Range("B7").Value = 9
Set-EmptyCell $extra.Range("C7")
Set-EmptyCell $extra.Range("D7")
Set-EmptyCell $extra.Range("E7")
Set-TextCell $extra.Range("F7") "NO"

# Row 8: MATCH_CODE 4566, batting position 10, 0 fours, 0 sixes, 1.14% runs
Set-TextCell $extra.Range("A8") "4566"
$extra.Range("B8").Value = 10
Set-TextCell $extra.Range("C8") "0"
Set-TextCell $extra.Range("D8") "0"
Set-TextCell $extra.Range("E8") "1.14%"
Set-TextCell $extra.Range("F8") "NO"

# Row 9: MATCH_CODE 4568, no batting position, no 4s/6s/pct, not MOM
Set-TextCell $extra.Range("A9") "4568"
Set-EmptyCell $extra.Range("B9")
Set-EmptyCell $extra.Range("C9")
Set-EmptyCell $extra.Range("D9")
Set-EmptyCell $extra.Range("E9")
Set-TextCell $extra.Range("F9") "NO"

$extra.Range("A1").Select()
